$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.045.89"
$ws.Range("E2").Value = "  +3.57%  "
$ws.Range("D3").Value = "3.451.79"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.21"
$ws.Range("E5").Value = "  +5.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.63"
$ws.Range("E6").Value = "  +7.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "3.445.59"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.15"
$ws.Range("E12").Value = "  +5.06%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.38"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "3.998.81"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "3.448.26"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "67.004.48"
$ws.Range("E18").Value = "  +3.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.08"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.02"
$ws.Range("E21").Value = "  +3.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "486.78"
$ws.Range("E22").Value = "  +8.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.33"
$ws.Range("E23").Value = "  +7.99%  "
$ws.Range("E24").Value = "  +23.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.41"
$ws.Range("E25").Value = "  +8.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.53"
$ws.Range("E26").Value = "  +3.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.95"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("E29").Value = "  +4.87%  "
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.25"
$ws.Range("E31").Value = "  +11.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "599.58"
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.73"
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.75"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  +7.11%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.66"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.385"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("D41").Value = "3.258.12"
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  +5.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0429"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("E45").Value = "  +24.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.28"
$ws.Range("E49").Value = "  +13.61%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.72"
$ws.Range("E50").Value = "  +5.73%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.12%  "
